# Create new graphic for scaling-alarms-capacity
# Reposition/resize three shapes on slide 1 to match the new artwork layout.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Group 243" -> shift left only (Top/Width/Height unchanged)
$grp243 = $s.Shapes.Item(1)
$grp243.Left = 222.0

# "Group 237" -> shift left only (Top/Width/Height unchanged)
$grp237 = $s.Shapes.Item(2)
$grp237.Left = 210.0

# "Picture 257" -> reposition and resize (rotated picture)
$pic257 = $s.Shapes.Item(23)
$pic257.Left = 209.3557
$pic257.Top = 136.60394
$pic257.Width = 37.4457
$pic257.Height = 37.7468
